# Macroferia Regional de Talca - Arveja Verde
# Insert a new data row at row 62 (pushing existing rows 62-119 down to 63-120)
# and populate it with a new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("62:62").Insert()

$ws.Range("A62").Value = 5
$ws.Range("B62").Value = "Macroferia Regional de Talca"
$ws.Range("C62").Value = "Maule"
$ws.Range("D62").Value = 44893
$ws.Range("E62").Value = 7
$ws.Range("F62").Value = 100112022
$ws.Range("G62").Value = "Arveja Verde"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 500
$ws.Range("K62").Value = 21000
$ws.Range("L62").Value = 21000
$ws.Range("M62").Value = 21000
$ws.Range("N62").Value = "`$/saco 25 kilos"
$ws.Range("O62").Value = "Región del Maule"
$ws.Range("P62").Value = 840
$ws.Range("Q62").Value = 25
$ws.Range("R62").Value = "Hortaliza"
